# T2237_Companies ... "Companies Changes 1 - 29th June 2023"
#
# The real change behind the resave is:
#   - Users sheet, A2: "Nicole Bicho" -> "Drew Koecher"
#   - Company sheet, N1's cell format loses its (redundant) explicit
#     number-format flag, ending up identical to the other bold header
#     cells (N1:M1) - i.e. just bold, default/General number format.
#
# Everything else in the raw XML diff (fileVersion/rupBuild bump, the
# xr/xr2/xr6/xr10 revision-tracking namespaces, the absPath username, the
# workbookView window geometry, etc.) is metadata Excel itself rewrites
# whenever a newer build resaves the file, not something driven by user
# actions in the object model, so it's not reproduced here.

$wb = $excel.ActiveWorkbook

# --- Company sheet: tidy up the "Postal Code" header cell's style so it
#     matches the other header cells (bold, default number format) ------
$company = $wb.Worksheets.Item("Company")
$company.Range("N1").Font.Bold = $true

# --- Users sheet: update the user name used by this duplicate-company
#     test data -----------------------------------------------------------
$users = $wb.Worksheets.Item("Users")
$users.Range("A2").Value = "Drew Koecher"
